# Data update: append the latest "Харьков" scrape results.
#  - Statistic sheet: new row 13 with the raw metrics + parsing status.
#  - Source sheet:    new row 29 with the rendered Ukrainian summary message.

$wb = $excel.ActiveWorkbook

# --- Statistic sheet: append row 13 -----------------------------------
$statistic = $wb.Worksheets.Item("Statistic")

$statistic.Cells.Item(13, 1).Value  = "2025-03-22 17:00:10.000000"
$statistic.Cells.Item(13, 2).Value  = 34905
$statistic.Cells.Item(13, 3).Value  = 61
$statistic.Cells.Item(13, 4).Value  = 382
$statistic.Cells.Item(13, 5).Value  = 1000
$statistic.Cells.Item(13, 6).Value  = 1700
$statistic.Cells.Item(13, 7).Value  = 9
$statistic.Cells.Item(13, 8).Value  = 3
$statistic.Cells.Item(13, 9).Value  = -4
$statistic.Cells.Item(13, 10).Value = "OK"

# --- Source sheet: append row 29 ---------------------------------------
$source = $wb.Worksheets.Item("Source")

$source.Cells.Item(29, 1).Value = "2025-03-22 17:00:10.000000"
$source.Cells.Item(29, 2).Value = 34905
$source.Cells.Item(29, 3).Value = "statistic"

$summary = "Статистика на Djinni за запитом **Харьков**:`nВакансій за 30 днів: 61`nКандидати онлайн: 382`nВилка по зарплаті: `$1000-1700`nВідгуків на одну вакансію: 9.0`nhttps://djinni.co/jobs/tg_search?keywords=%D0%A5%D0%B0%D1%80%D1%8C%D0%BA%D0%BE%D0%B2`nАктивність за тиждень:`nВакансій: +3`nКандидатів: -4"

$source.Cells.Item(29, 4).Value = $summary

# Setting a multi-line value auto-expands the row height (customHeight flag);
# AutoFit puts the row back to a content-derived, non-pinned height so the
# row element stays bare, matching the other data rows on this sheet.
$source.Rows.Item(29).EntireRow.AutoFit()
